$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of LinkedIn/article content scraped after the last existing entry (rows 322-336).
# Data: Title (A), Summary (B), URL (C), Date (D) - stored as plain text, same as existing rows.
$newRows = @(
    ,@('At Copeland, we’re proud to advance sustainable, efficient, and safe production systems for compressors and heat pumps. Our', 'At Copeland, we’re proud to advance sustainable, efficient, and safe production systems for compressors and heat pumps. Our latest interview with Ong-art Veerachartyanukul, General Manager of our...', 'https://www.linkedin.com/feed/update/urn:li:activity:7334126665112924160', '2025-05-30')
    ,@(' SRMTec in 60 seconds: technology, precision, passion. We manufacture screw and piston compressors for refrigeration', ' SRMTec in 60 seconds: technology, precision, passion.We manufacture screw and piston compressors for refrigeration, fully designed and made in Italy. Each compressor is tested twice: performance and...', 'https://www.linkedin.com/feed/update/urn:li:activity:7334111489491881984', '2025-05-30')
    ,@('Our innovative technologies are driving the shift toward sustainable cooling by promoting low-energy consumption and environment', 'Our innovative technologies are driving the shift toward sustainable cooling by promoting low-energy consumption and environmentally responsible refrigerants.From cutting-edge compressor designs to...', 'https://www.linkedin.com/feed/update/urn:li:activity:7334073630852583424', '2025-05-30')
    ,@('We are excited to announce Andy Baker, MBA as Copeland’s Vice President of Commercial Excellence, a new role within Copeland tha', 'We are excited to announce Andy Baker, MBA as Copeland’s Vice President of Commercial Excellence, a new role within Copeland that will be responsible for leading the global implementation of...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333915424494952450', '2025-05-29')
    ,@('Industrial applications demand reliability, precision and efficiency. Copeland’s innovative Vilter single-screw compression', 'Industrial applications demand reliability, precision and efficiency. Copeland’s innovative Vilter single-screw compression technology has proven to be up to the task, delivering superior performance...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333885283546943490', '2025-05-29')
    ,@('We were thrilled to host the North Florida Refrigerating Engineers & Technicians Association (RETA) - Official Page Chapter meet', 'We were thrilled to host the North Florida Refrigerating Engineers & Technicians Association (RETA) - Official Page Chapter meeting at our #FrascoldUSA office, the first chapter meeting of the...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333786437319819264', '2025-05-29')
    ,@('Mark your calendars for June 12, 2025, at 11 am CET, as we’re hosting an online session – Retrofitting Commercial Buildings', 'Mark your calendars for June 12, 2025, at 11 am CET, as we’re hosting an online session – Retrofitting Commercial Buildings with Modulating Room Control. Discover how transitioning from traditional...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333779205492523008', '2025-05-29')
    ,@('Copeland has partnered with Ecozen to deliver solar-powered cold room systems that enable efficient, off-grid refrigeration for ', 'Copeland has partnered with Ecozen to deliver solar-powered cold room systems that enable efficient, off-grid refrigeration for the agricultural sector in India. Powered by our variable speed...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333711238180327426', '2025-05-29')
    ,@(' Join Our Dynamic Team at Snowman Group!  We are excited to announce the following openings:  Sale', ' Join Our Dynamic Team at Snowman Group! We are excited to announce the following openings: Sales Engineer (Food chain industry; Chemical, Oil & Gas industry) After-sales Service EngineerIf you’re...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333677525769895936', '2025-05-29')
    ,@('Copeland Sensi Equipment Interface Module', 'Copeland’s Sensi Equipment Interface Module (Indoor + Outdoor) is a new interface module that eliminates the need for pulling additional HVAC wiring to the thermostat, helping solve installation...', 'https://www.phcppros.com/articles/21585-copeland-sensi-equipment-interface-module', '2025-05-28')
    ,@('Copeland Thailand Innovation Drives Sustainability, Efficiency, & Safety Wins', ' JARN interviewed Ong-Art Veerachartyanukul, general manager of Copeland’s facility in Rayong, Thailand, who talked about Copeland’s advancing', 'https://www.ejarn.com/article/detail/88911', '2025-05-28')
    ,@('Innovations in refrigeration are shaping a more sustainable, efficient future across North America. On June 11 at the ATMOsphere', 'Innovations in refrigeration are shaping a more sustainable, efficient future across North America. On June 11 at the ATMOsphere America Summit 2025 in Atlanta, Andre Patenaude, CET, director of...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333477479543672832', '2025-05-28')
    ,@('NeoCharge is transforming industrial refrigeration by reducing refrigerant charge in both new and existing systems. With easy in', 'NeoCharge is transforming industrial refrigeration by reducing refrigerant charge in both new and existing systems. With easy installation, it offers significant energy savings and increased capacity...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333454503226077184', '2025-05-28')
    ,@('From Deep-Sea Fishing to End-Point RetailSnowman Refrigeration Technology Locks in Every Bite of Freshness. See you at CHI', 'From Deep-Sea Fishing to End-Point RetailSnowman Refrigeration Technology Locks in Every Bite of Freshness.See you at CHINA (FUZHOU)INTERNATIONAL SEAFOOD & FISHERIES EXPO!  Venue: Fuzhou Strait...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333431848183730176', '2025-05-28')
    ,@(' Looking back on an inspiring day: “#Hydrocarbons in the Future of #HVACR” – Padua, May 23. The energ', ' Looking back on an inspiring day: “#Hydrocarbons in the Future of #HVACR” – Padua, May 23.The energy, insights and meaningful exchanges of May 23 are still with us.A day that brought together...', 'https://www.linkedin.com/feed/update/urn:li:activity:7333415061496958977', '2025-05-28')
)

$startRow = 322
$endRow = $startRow + $newRows.Count - 1

# Force the Date column to Text format first so Excel does not auto-convert
# the "yyyy-mm-dd" strings into date serial numbers (existing rows store dates as text).
$ws.Range("D{0}:D{1}" -f $startRow, $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}

Write-Host ("New used range: " + $ws.UsedRange.Address())
